# Append the new "suite tp banque avec repo et view" entries after the
# last existing paragraph ("Créer un service qui permet de hasher et de
# vérifier le hash d'un mot de passe"), mirroring the formatting of the
# surrounding paragraphs (Helvetica, 12pt).

$d = $word.ActiveDocument

$newParagraphs = @(
    "",
    "",
    "List nugets",
    "",
    "EntityFrameWorkcore",
    "EntityFrameWorkCore tools",
    "EntityFrameWork"
)

$anchor = $d.Paragraphs.Last.Range

foreach ($text in $newParagraphs) {
    $anchor.InsertParagraphAfter()
    $anchor = $d.Paragraphs.Last.Range
    $anchor.Text = $text
}
